# #129 adding support of named ranges for DataModel class
#
# B2 previously concatenated two cell references (A1, A2). The edit swaps
# the first argument for a bare, undefined name ("text") so the formula now
# references a (currently missing) named range instead of cell A1 — which is
# exactly the scenario this "Error_Name" workbook is meant to exercise: the
# formula evaluates to a #NAME? error.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Formula = "=CONCATENATE(text,A2)"
